$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Rules")

# Restore value: cell C10 (row 10) changes from 18 to 1.
$ws.Range("C10").Value = 1
